# Add a new 'alt' column before the existing SC_58SM column on the
# "variant comparison" sheet. Inserting the column shifts I:N -> J:O
# (including the header cells and all eight data rows), matching the
# target workbook's new A1:O8 dimension.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variant comparison")

# Shift columns I:N one place to the right, starting a new column I.
$ws.Columns.Item(9).Insert()

# New column header.
$ws.Range("I1").Value = "alt"

# New column's per-row values: matches the sample's existing "ref"/allele
# style columns next to it (only row 2 carries a value; the rest are blank
# inline strings, same as the other sample columns in those rows).
$ws.Range("I2").Value = "deletion"
$ws.Range("I3").Value = ""
$ws.Range("I4").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("I8").Value = ""
